$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, shifting existing rows 232:255 down to 233:256
$ws.Rows("232:232").Insert()

# Populate the newly inserted row 232 with the new record
$ws.Range("A232").Value = 1
$ws.Range("B232").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C232").Value = "Arica y Parinacota"
$ws.Range("D232").Value = 44748
$ws.Range("E232").Value = 15
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100102
$ws.Range("H232").Value = "Cítricos"
$ws.Range("I232").Value = 100102003
$ws.Range("J232").Value = "Limón"
$ws.Range("K232").Value = "Sin especificar"
$ws.Range("L232").Value = "2a amarillo"
$ws.Range("M232").Value = 300
$ws.Range("N232").Value = 9000
$ws.Range("O232").Value = 10000
$ws.Range("P232").Value = 9500
$ws.Range("Q232").Value = "$/caja 20 kilos"
$ws.Range("R232").Value = "Región de Coquimbo"
$ws.Range("S232").Value = 475
$ws.Range("T232").Value = 20
